$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend dimension implicitly by writing data; set number format for date-like text columns first
$ws.Range("S221:S229").NumberFormat = "@"

# Row 221
$ws.Range("B221").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("C221").Value = 2026
$ws.Range("D221").Value = "MDPI AG"
$ws.Range("E221").Value = "Energies"
$ws.Range("F221").Value = "Khan, Waqar Hussain; Ahn, Dukju"
$ws.Range("H221").Value = "10.3390/en19041063"
$ws.Range("I221").Value = "https://doi.org/10.3390/en19041063"
$ws.Range("J221").Value = "Journal"
$ws.Range("K221").Value = "Inverter"
$ws.Range("L221").Value = "Experiment"
$ws.Range("M221").Value = "Contacts"
$ws.Range("Q221").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("R221").Value = "High"
$ws.Range("S221").Value = "2026-02-19"

# Row 222
$ws.Range("B222").Value = "Distilling Protein Language Models with Complementary Regularizers"
$ws.Range("C222").Value = 2026
$ws.Range("D222").Value = "openRxiv"
$ws.Range("F222").Value = "Wijaya, Edward"
$ws.Range("H222").Value = "10.64898/2026.02.17.706304"
$ws.Range("I222").Value = "https://doi.org/10.64898/2026.02.17.706304"
$ws.Range("J222").Value = "Journal"
$ws.Range("K222").Value = "Co-integration"
$ws.Range("L222").Value = "Experiment"
$ws.Range("M222").Value = "Gate Stack"
$ws.Range("Q222").Value = "Distilling Protein Language Models with Complementary Regularizers"
$ws.Range("R222").Value = "High"
$ws.Range("S222").Value = "2026-02-19"

# Row 223
$ws.Range("B223").Value = "Distilling Protein Language Models with Complementary Regularizers"
$ws.Range("C223").Value = 2026
$ws.Range("D223").Value = "openRxiv"
$ws.Range("F223").Value = "Wijaya, Edward"
$ws.Range("H223").Value = "10.64898/2026.02.17.706304"
$ws.Range("I223").Value = "https://doi.org/10.64898/2026.02.17.706304"
$ws.Range("J223").Value = "Journal"
$ws.Range("K223").Value = "Co-integration"
$ws.Range("L223").Value = "Experiment"
$ws.Range("M223").Value = "Gate Stack"
$ws.Range("Q223").Value = "Distilling Protein Language Models with Complementary Regularizers"
$ws.Range("R223").Value = "High"
$ws.Range("S223").Value = "2026-02-19"

# Row 224
$ws.Range("B224").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("C224").Value = 2026
$ws.Range("D224").Value = "AIP Publishing"
$ws.Range("E224").Value = "Journal of Applied Physics"
$ws.Range("F224").Value = "Asteris, Aias; Nguyen, Thai-Son; Chang, Chuan F. C.; Savant, Chandrashekhar; Lonergan, Pierce; Xing, Huili G.; Jena, Debdeep"
$ws.Range("H224").Value = "10.1063/5.0312252"
$ws.Range("I224").Value = "https://doi.org/10.1063/5.0312252"
$ws.Range("J224").Value = "Journal"
$ws.Range("K224").Value = "n-FET"
$ws.Range("L224").Value = "Experiment"
$ws.Range("M224").Value = "Gate Stack"
$ws.Range("Q224").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("R224").Value = "High"
$ws.Range("S224").Value = "2026-02-19"

# Row 225
$ws.Range("B225").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("C225").Value = 2026
$ws.Range("D225").Value = "AIP Publishing"
$ws.Range("E225").Value = "Journal of Applied Physics"
$ws.Range("F225").Value = "Asteris, Aias; Nguyen, Thai-Son; Chang, Chuan F. C.; Savant, Chandrashekhar; Lonergan, Pierce; Xing, Huili G.; Jena, Debdeep"
$ws.Range("H225").Value = "10.1063/5.0312252"
$ws.Range("I225").Value = "https://doi.org/10.1063/5.0312252"
$ws.Range("J225").Value = "Journal"
$ws.Range("K225").Value = "n-FET"
$ws.Range("L225").Value = "Experiment"
$ws.Range("M225").Value = "Gate Stack"
$ws.Range("Q225").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("R225").Value = "High"
$ws.Range("S225").Value = "2026-02-19"

# Row 226
$ws.Range("B226").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("C226").Value = 2026
$ws.Range("D226").Value = "MDPI AG"
$ws.Range("E226").Value = "Energies"
$ws.Range("F226").Value = "Khan, Waqar Hussain; Ahn, Dukju"
$ws.Range("H226").Value = "10.3390/en19041063"
$ws.Range("I226").Value = "https://doi.org/10.3390/en19041063"
$ws.Range("J226").Value = "Journal"
$ws.Range("K226").Value = "Inverter"
$ws.Range("L226").Value = "Experiment"
$ws.Range("M226").Value = "Contacts"
$ws.Range("Q226").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("R226").Value = "High"
$ws.Range("S226").Value = "2026-02-19"

# Row 227
$ws.Range("B227").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("C227").Value = 2026
$ws.Range("D227").Value = "MDPI AG"
$ws.Range("E227").Value = "Energies"
$ws.Range("F227").Value = "Khan, Waqar Hussain; Ahn, Dukju"
$ws.Range("H227").Value = "10.3390/en19041063"
$ws.Range("I227").Value = "https://doi.org/10.3390/en19041063"
$ws.Range("J227").Value = "Journal"
$ws.Range("K227").Value = "Inverter"
$ws.Range("L227").Value = "Experiment"
$ws.Range("M227").Value = "Contacts"
$ws.Range("Q227").Value = "Current-Summing Multilevel LCC Inverter for Radiated EMI Harmonic Reduction in Wireless Power Transfer"
$ws.Range("R227").Value = "High"
$ws.Range("S227").Value = "2026-02-19"

# Row 228
$ws.Range("B228").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("C228").Value = 2026
$ws.Range("D228").Value = "AIP Publishing"
$ws.Range("E228").Value = "Journal of Applied Physics"
$ws.Range("F228").Value = "Asteris, Aias; Nguyen, Thai-Son; Chang, Chuan F. C.; Savant, Chandrashekhar; Lonergan, Pierce; Xing, Huili G.; Jena, Debdeep"
$ws.Range("H228").Value = "10.1063/5.0312252"
$ws.Range("I228").Value = "https://doi.org/10.1063/5.0312252"
$ws.Range("J228").Value = "Journal"
$ws.Range("K228").Value = "n-FET"
$ws.Range("L228").Value = "Experiment"
$ws.Range("M228").Value = "Gate Stack"
$ws.Range("Q228").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("R228").Value = "High"
$ws.Range("S228").Value = "2026-02-19"

# Row 229
$ws.Range("B229").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("C229").Value = 2026
$ws.Range("D229").Value = "AIP Publishing"
$ws.Range("E229").Value = "Journal of Applied Physics"
$ws.Range("F229").Value = "Asteris, Aias; Nguyen, Thai-Son; Chang, Chuan F. C.; Savant, Chandrashekhar; Lonergan, Pierce; Xing, Huili G.; Jena, Debdeep"
$ws.Range("H229").Value = "10.1063/5.0312252"
$ws.Range("I229").Value = "https://doi.org/10.1063/5.0312252"
$ws.Range("J229").Value = "Journal"
$ws.Range("K229").Value = "n-FET"
$ws.Range("L229").Value = "Experiment"
$ws.Range("M229").Value = "Gate Stack"
$ws.Range("Q229").Value = "High mobility multiple-channel AlScN/GaN heterostructures"
$ws.Range("R229").Value = "High"
$ws.Range("S229").Value = "2026-02-19"
